# V9.0 Wprowadzanie swojego nicku
#
# Adds a "type your own nickname" row group: six new entries land at the
# top of the data (rows 2-7), pushing the existing rows down by six and
# renumbering them (column A) and updating some of their figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the current data rows (2-20) down to rows 8-26, carrying their
# content/format with them.
$ws.Rows("2:7").Insert()

# Insert() clones the header row's formatting onto the freshly inserted
# rows; strip that back off so they start out unstyled like every other
# data row (only column A keeps the "s=1" style, applied below).
$ws.Range("B2:E7").ClearFormats()

# Column A uses the same bold/bordered/centered style as the rest of the
# data rows (cellXf index 1) - clone it from an already-styled cell below
# rather than rebuilding it property-by-property, which would otherwise
# mint a near-duplicate style.
$ws.Cells.Item(8, 1).Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---- New rows 2-7: custom-nickname entries ----
$newRows = @(
    @(0, "", "Extreme", "Dolnośląskie", 30),
    @(1, "maziar", "Hard", "Podlaskie", 8),
    @(2, "Maksssssssssss", "Hard", "Łódzkie", 8),
    @(3, "Maks", "Hard", "Podlaskie", 10),
    @(4, "", "Extreme", "Podlaskie", 12),
    @(5, "", "Extreme", "Dolnośląskie", 36)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

# ---- Shifted rows, now 8-26: refresh to the updated figures ----
$shiftedRows = @(
    @(6, "maks", "Extreme", "Dolnośląskie", 21),
    @(7, "maks", "Hard", "Dolnośląskie", 16),
    @(8, "maks", "Hard", "Podlaskie", 16),
    @(9, "maks", "Extreme", "Wszystkie", 485),
    @(10, "maks", "Extreme", "Wszystkie", 540),
    @(11, "maks", "Extreme", "Pomorskie", 15),
    @(12, "maks", "Hard", "Pomorskie", 8),
    @(13, "maks", "Medium", "Pomorskie", 2),
    @(14, "maks", "Medium", "Pomorskie", 5),
    @(15, "maks", "Medium", "Pomorskie", 4),
    @(16, "maks", "Medium", "Pomorskie", 5),
    @(17, "maks", "Medium", "Pomorskie", 5),
    @(18, "maks", "Easy", "Wielkopolskie", 13),
    @(19, "maks", "Extreme", "Lubuskie", 9),
    @(20, "maks", "Extreme", "Opolskie", 6),
    @(21, "maks", "Extreme", "Śląskie", 33),
    @(22, "maks", "Extreme", "Śląskie", 33),
    @(23, "maks", "Extreme", "Kujawsko-Pomorskie", 18),
    @(24, "maks", "Extreme", "Kujawsko-Pomorskie", 18)
)

for ($i = 0; $i -lt $shiftedRows.Count; $i++) {
    $r = 8 + $i
    $vals = $shiftedRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

$ws.Range("A1").Select()
